# Plau-Igf2r NATMI output update ("Natmi following Dr Hou advice").
#
# The ligand-/receptor-expressing cell counts for every row go from 1 to 3
# (columns E and K), which changes every downstream value that is derived
# from those counts: total-expression columns (H, N, R), the specificity
# columns that are re-normalised against the new per-(sending-cluster)
# totals (I/J, O/P), and the edge-weight columns that multiply averages and
# totals together (Q, R, S, T). Column G (ligand average expression) and M
# (receptor average expression) also move slightly because they come from
# the same underlying recompute. F, L (detection rates) and the text
# columns A-D are unchanged.
#
# All seventeen data rows (r=2..17) keep the same row/column layout, so we
# just overwrite the affected cells in each row with their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 19.434592
$ws.Range("H2").Value = 58.303776
$ws.Range("I2").Value = 0.1244167820899015
$ws.Range("J2").Value = 0.1244167820899015
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.81388166666666
$ws.Range("N2").Value = 92.441645
$ws.Range("O2").Value = 0.1674303154124352
$ws.Range("P2").Value = 0.1674303154124352
$ws.Range("Q2").Value = 598.8552181279466
$ws.Range("R2").Value = 5389.696963151519
$ws.Range("S2").Value = 0.02083114106791243
$ws.Range("T2").Value = 0.02083114106791243

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 19.434592
$ws.Range("H3").Value = 58.303776
$ws.Range("I3").Value = 0.1244167820899015
$ws.Range("J3").Value = 0.1244167820899015
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 82.73043066666666
$ws.Range("N3").Value = 248.191292
$ws.Range("O3").Value = 0.4495240895180934
$ws.Range("P3").Value = 0.4495240895180934
$ws.Range("Q3").Value = 1607.832165990955
$ws.Range("R3").Value = 14470.48949391859
$ws.Range("S3").Value = 0.05592834068973401
$ws.Range("T3").Value = 0.05592834068973401

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 19.434592
$ws.Range("H4").Value = 58.303776
$ws.Range("I4").Value = 0.1244167820899015
$ws.Range("J4").Value = 0.1244167820899015
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.871077999999999
$ws.Range("N4").Value = 26.613234
$ws.Range("O4").Value = 0.04820189172060866
$ws.Range("P4").Value = 0.04820189172060866
$ws.Range("Q4").Value = 172.405781530176
$ws.Range("R4").Value = 1551.652033771584
$ws.Range("S4").Value = 0.005997124258523997
$ws.Range("T4").Value = 0.005997124258523996

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.434592
$ws.Range("H5").Value = 58.303776
$ws.Range("I5").Value = 0.1244167820899015
$ws.Range("J5").Value = 0.1244167820899015
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 61.624648
$ws.Range("N5").Value = 184.873944
$ws.Range("O5").Value = 0.3348437033488628
$ws.Range("P5").Value = 0.3348437033488628
$ws.Range("Q5").Value = 1197.649891023616
$ws.Range("R5").Value = 10778.84901921254
$ws.Range("S5").Value = 0.0416601760737311
$ws.Range("T5").Value = 0.04166017607373109

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 48.891945
$ws.Range("H6").Value = 146.675835
$ws.Range("I6").Value = 0.3129974875220664
$ws.Range("J6").Value = 0.3129974875220664
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.81388166666666
$ws.Range("N6").Value = 92.441645
$ws.Range("O6").Value = 0.1674303154124352
$ws.Range("P6").Value = 0.1674303154124352
$ws.Range("Q6").Value = 1506.550607683175
$ws.Range("R6").Value = 13558.95546914858
$ws.Range("S6").Value = 0.05240526805911934
$ws.Range("T6").Value = 0.05240526805911933

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 48.891945
$ws.Range("H7").Value = 146.675835
$ws.Range("I7").Value = 0.3129974875220664
$ws.Range("J7").Value = 0.3129974875220664
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 82.73043066666666
$ws.Range("N7").Value = 248.191292
$ws.Range("O7").Value = 0.4495240895180934
$ws.Range("P7").Value = 0.4495240895180934
$ws.Range("Q7").Value = 4044.85166598098
$ws.Range("R7").Value = 36403.66499382882
$ws.Range("S7").Value = 0.1406999105998077
$ws.Range("T7").Value = 0.1406999105998077

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 48.891945
$ws.Range("H8").Value = 146.675835
$ws.Range("I8").Value = 0.3129974875220664
$ws.Range("J8").Value = 0.3129974875220664
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.871077999999999
$ws.Range("N8").Value = 26.613234
$ws.Range("O8").Value = 0.04820189172060866
$ws.Range("P8").Value = 0.04820189172060866
$ws.Range("Q8").Value = 433.7242576667099
$ws.Range("R8").Value = 3903.51831900039
$ws.Range("S8").Value = 0.01508707100236121
$ws.Range("T8").Value = 0.0150870710023612

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 48.891945
$ws.Range("H9").Value = 146.675835
$ws.Range("I9").Value = 0.3129974875220664
$ws.Range("J9").Value = 0.3129974875220664
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 61.624648
$ws.Range("N9").Value = 184.873944
$ws.Range("O9").Value = 0.3348437033488628
$ws.Range("P9").Value = 0.3348437033488628
$ws.Range("Q9").Value = 3012.94890066036
$ws.Range("R9").Value = 27116.54010594324
$ws.Range("S9").Value = 0.1048052378607782
$ws.Range("T9").Value = 0.1048052378607782

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 65.19353
$ws.Range("H10").Value = 195.58059
$ws.Range("I10").Value = 0.4173573191390618
$ws.Range("J10").Value = 0.4173573191390618
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.81388166666666
$ws.Range("N10").Value = 92.441645
$ws.Range("O10").Value = 0.1674303154124352
$ws.Range("P10").Value = 0.1674303154124352
$ws.Range("Q10").Value = 2008.865718852283
$ws.Range("R10").Value = 18079.79146967055
$ws.Range("S10").Value = 0.0698782675831415
$ws.Range("T10").Value = 0.0698782675831415

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 65.19353
$ws.Range("H11").Value = 195.58059
$ws.Range("I11").Value = 0.4173573191390618
$ws.Range("J11").Value = 0.4173573191390618
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 82.73043066666666
$ws.Range("N11").Value = 248.191292
$ws.Range("O11").Value = 0.4495240895180934
$ws.Range("P11").Value = 0.4495240895180934
$ws.Range("Q11").Value = 5393.488813580253
$ws.Range("R11").Value = 48541.39932222228
$ws.Range("S11").Value = 0.1876121688896991
$ws.Range("T11").Value = 0.1876121688896991

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 65.19353
$ws.Range("H12").Value = 195.58059
$ws.Range("I12").Value = 0.4173573191390618
$ws.Range("J12").Value = 0.4173573191390618
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 8.871077999999999
$ws.Range("N12").Value = 26.613234
$ws.Range("O12").Value = 0.04820189172060866
$ws.Range("P12").Value = 0.04820189172060866
$ws.Range("Q12").Value = 578.3368897253399
$ws.Range("R12").Value = 5205.03200752806
$ws.Range("S12").Value = 0.02011741230594457
$ws.Range("T12").Value = 0.02011741230594457

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 65.19353
$ws.Range("H13").Value = 195.58059
$ws.Range("I13").Value = 0.4173573191390618
$ws.Range("J13").Value = 0.4173573191390618
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 61.624648
$ws.Range("N13").Value = 184.873944
$ws.Range("O13").Value = 0.3348437033488628
$ws.Range("P13").Value = 0.3348437033488628
$ws.Range("Q13").Value = 4017.52833812744
$ws.Range("R13").Value = 36157.75504314696
$ws.Range("S13").Value = 0.1397494703602767
$ws.Range("T13").Value = 0.1397494703602767

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 22.685484
$ws.Range("H14").Value = 68.05645200000001
$ws.Range("I14").Value = 0.1452284112489703
$ws.Range("J14").Value = 0.1452284112489703
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.81388166666666
$ws.Range("N14").Value = 92.441645
$ws.Range("O14").Value = 0.1674303154124352
$ws.Range("P14").Value = 0.1674303154124352
$ws.Range("Q14").Value = 699.02781952706
$ws.Range("R14").Value = 6291.250375743541
$ws.Range("S14").Value = 0.02431563870226195
$ws.Range("T14").Value = 0.02431563870226195

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 22.685484
$ws.Range("H15").Value = 68.05645200000001
$ws.Range("I15").Value = 0.1452284112489703
$ws.Range("J15").Value = 0.1452284112489703
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 82.73043066666666
$ws.Range("N15").Value = 248.191292
$ws.Range("O15").Value = 0.4495240895180934
$ws.Range("P15").Value = 0.4495240895180934
$ws.Range("Q15").Value = 1876.779861201776
$ws.Range("R15").Value = 16891.01875081599
$ws.Range("S15").Value = 0.0652836693388526
$ws.Range("T15").Value = 0.0652836693388526

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 22.685484
$ws.Range("H16").Value = 68.05645200000001
$ws.Range("I16").Value = 0.1452284112489703
$ws.Range("J16").Value = 0.1452284112489703
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 8.871077999999999
$ws.Range("N16").Value = 26.613234
$ws.Range("O16").Value = 0.04820189172060866
$ws.Range("P16").Value = 0.04820189172060866
$ws.Range("Q16").Value = 201.244698031752
$ws.Range("R16").Value = 1811.202282285768
$ws.Range("S16").Value = 0.007000284153778891
$ws.Range("T16").Value = 0.007000284153778891

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 22.685484
$ws.Range("H17").Value = 68.05645200000001
$ws.Range("I17").Value = 0.1452284112489703
$ws.Range("J17").Value = 0.1452284112489703
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 61.624648
$ws.Range("N17").Value = 184.873944
$ws.Range("O17").Value = 0.3348437033488628
$ws.Range("P17").Value = 0.3348437033488628
$ws.Range("Q17").Value = 1397.984966209632
$ws.Range("R17").Value = 12581.86469588669
$ws.Range("S17").Value = 0.04862881905407686
$ws.Range("T17").Value = 0.04862881905407686

